# "Generate Report for Handoff" — refresh the localization-status report:
#   - Status flips from "Handed back: in sync with en-US" to "Ready for handoff"
#   - Latest Handoff / HO Xliff Generate timestamps bump forward a bit
#   - Status/Datetime columns on each sheet get narrower (report layout tweak)

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---------------------------------------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-23 00:56:17"

# --- zh-cn sheet --------------------------------------------------------
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-23 00:56:12"

# --- de-de sheet --------------------------------------------------------
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-23 00:56:17"

# --- Narrower Status / Datetime columns on every sheet ------------------
# (ColumnWidth is in "characters"; the host snaps to a 1/6-character pixel
#  grid on write, so feed it the pre-image that lands on the nearest grid
#  point to the target ~17.22-character width.)
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333
$zhcn.Columns.Item(3).ColumnWidth = 16.3333333333333
$dede.Columns.Item(3).ColumnWidth = 16.3333333333333
